# Edit: reorder "Periodo Mora" rows (16-39) from ascending to descending
# and keep the "Valor Mora" figures attached to their correct period
# (period 2005 keeps 19791, all other periods keep 31249).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Periods in descending order (2005 .. 1806), mirroring the previously
# ascending order (1806 .. 2005) that occupied rows 16-39.
$periods = @(
    "2005","2004","2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# Valor Mora (column F) for the first row (now period 2005) is the
# partial-period value 19791; the last row (now period 1806) reverts to
# the standard 31249 value (swap of the previous F16/F39 contents).
$ws.Cells.Item(16, 6).Value = 19791
$ws.Cells.Item(39, 6).Value = 31249
